$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Kate_hours")
$ws2 = $wb.Worksheets.Item("Ben_hours")

# --- Add new logged-hours rows to Kate_hours (rows 9-13) ---
$ws1.Range("A9").Value = 45029
$ws1.Range("A9").NumberFormat = "d-mmm"
$ws1.Range("B9").Value = 2
$ws1.Range("C9").Value = "set up worms records search, taxonomy"

$ws1.Range("A10").Value = 45030
$ws1.Range("A10").NumberFormat = "d-mmm"
$ws1.Range("B10").Value = 2
$ws1.Range("C10").Value = "comparison lists, region setup transferred"

$ws1.Range("A11").Value = 45034
$ws1.Range("A11").NumberFormat = "d-mmm"
$ws1.Range("B11").Value = 0.5
$ws1.Range("C11").Value = "documentation, helping"

$ws1.Range("A12").Value = 45036
$ws1.Range("A12").NumberFormat = "d-mmm"
$ws1.Range("B12").Value = 1
$ws1.Range("C12").Value = "meeting about goals"

$ws1.Range("A13").Value = 45037
$ws1.Range("A13").NumberFormat = "d-mmm"
$ws1.Range("B13").Value = 1
$ws1.Range("C13").Value = "workflow planning"

# --- Update active sheet / selection to match the author's saved view ---
[void]$ws1.Activate()
[void]$ws1.Range("H13").Select()
